# Regenerate the "K" (strikeouts) column (G) for save_data rows 2-13 using
# the corrected per-game values (was previously populated from a different
# "Strike#" source; now uses the K count directly).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 4
    3  = 4
    4  = 7
    5  = 6
    6  = 2
    7  = 4
    8  = 6
    9  = 5
    10 = 3
    11 = 5
    12 = 3
    13 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
